# Update beta coefficients column (B) to reflect the "30 day" re-estimation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = -0.0074257969549297353
$ws.Range("B4").Value  = -0.021164895635322683
$ws.Range("B5").Value  = -0.070302205579692134
$ws.Range("B6").Value  = -0.20254474077803125
$ws.Range("B7").Value  = 0.58163984142033165
$ws.Range("B8").Value  = -0.19014030198723705
$ws.Range("B10").Value = 0.12663625318142149
$ws.Range("B12").Value = -0.088178983128388566
$ws.Range("B13").Value = 0.098116295064896664
$ws.Range("B14").Value = -0.10164986187615825
